# Updated symbol list on Thu Dec 29 13:47:30 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores values as text (they look numeric, e.g.
# "245.73", but must stay text so things like trailing zeros - "6.500" -
# and leading zeros - "0.00005212" - round-trip exactly, matching how the
# source sheet keeps every cell as inline/text strings). Force the
# NumberFormat to Text ("@") before writing those so Excel doesn't
# reinterpret the literal as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($a1, $value) {
    $rng = $ws.Range($a1)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- Price (column D) updates ---
Set-TextValue "D2"  "245.84"
Set-TextValue "D3"  "24.19"
Set-TextValue "D4"  "5.301"
Set-TextValue "D5"  "0.05741"
Set-TextValue "D6"  "6.500"
Set-TextValue "D8"  "0.8174"
Set-TextValue "D9"  "0.8718"
Set-TextValue "D11" "0.06994"
Set-TextValue "D12" "0.03178"
Set-TextValue "D13" "0.02932"
Set-TextValue "D14" "0.09377"
Set-TextValue "D15" "3.730"
Set-TextValue "D16" "0.001522"
Set-TextValue "D17" "0.04723"
Set-TextValue "D18" "0.0006004"
Set-TextValue "D19" "0.006185"
Set-TextValue "D20" "0.001237"
Set-TextValue "D21" "0.003862"
Set-TextValue "D24" "2.150"
Set-TextValue "D26" "0.1331"
Set-TextValue "D28" "0.0003012"

# Row 28 (UpBots) "best in 24h" tag appended to its Data id.
$ws.Range("E28").Value = "27UpBotsUBXTBestin24h"

Set-TextValue "D40" "0.03721"

# Rows 41-43 got re-ranked: KickToken / BKEXToken / CEJI rotate down one
# row, with BKEXToken now in the row that used to be KickToken, CEJI in
# the row that used to be BKEXToken, and KickToken in the row that used
# to be CEJI (each with refreshed price + link + Data id).
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1059"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.002227"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.006394"
$ws.Range("E43").Value = "42KickTokenKICK"

Set-TextValue "D44" "0.007634"
Set-TextValue "D45" "0.00005224"
Set-TextValue "D47" "0.3597"
Set-TextValue "D48" "0.002712"
